# Core Module Description.docx -- "namespace string ed:3.0 -> ed-2"
#
# 1. Namespace paragraph: collapse the 3 split runs ("urn:" / "iso:std" /
#    ":iso:30042:ed:3.0", separated by a gramStart/gramEnd proofErr pair)
#    into a single run with the new text "urn:iso:std:iso:30042:ed-2", and
#    drop the grammar-check proofErr markers.
# 2. The _GoBack bookmark (Word's "last edit" marker) follows the edit, so
#    it moves from the "Core module ... valid TBX" paragraph to the end of
#    the Namespace paragraph.
# 3. In the "Core module..." paragraph, re-split the run around the word
#    "classificationElement" and mark it with a spellStart/spellEnd
#    proofErr pair (Word's spell checker doesn't know this camelCase XML
#    tag name).
# 4. Several data-category names in the table (transacGrp, conceptEntry,
#    langSec, termSec, adminGrp, termNoteGrp, titleStmt, transacGrp,
#    termSec) are likewise camelCase XML tag/attribute names that Word's
#    spell checker flags; wrap each in its own spellStart/spellEnd pair.

$d = $word.ActiveDocument

$pkgHeader = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml" pkg:padding="512"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">'
$pkgFooter = '</w:document></pkg:xmlData></pkg:part></pkg:package>'

function Set-RangeXml($range, [string]$bodyXml) {
    $range.InsertXML($pkgHeader + '<w:body>' + $bodyXml + '</w:body>' + $pkgFooter)
}

# ---------------------------------------------------------------------
# 1 & 2. Namespace paragraph text + bookmark relocation
# ---------------------------------------------------------------------

$nsPara = $d.Paragraphs(4).Range
$nsSub = $d.Range($nsPara.Start, $nsPara.End - 1)
Set-RangeXml $nsSub '<w:p><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Namespace</w:t></w:r><w:r><w:rPr><w:rStyle w:val="FootnoteReference"/><w:b/><w:bCs/></w:rPr><w:footnoteReference w:id="1"/></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">: </w:t></w:r><w:r><w:t>urn:iso:std:iso:30042:ed-2</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'

# ---------------------------------------------------------------------
# 3. "The Core module ..." paragraph: split out classificationElement,
#    drop the old bookmark location (superseded by the move above).
# ---------------------------------------------------------------------

$corePara = $d.Paragraphs(6).Range
$coreSub = $d.Range($corePara.Start, $corePara.End - 1)
$coreBody = '<w:p>' +
  '<w:r><w:tab/></w:r>' +
  '<w:r><w:t>The Core module is the foundational module for all valid TBX dialects.  Unlike the other modules, Core module provides the XML structure, elements, and attributes which define TBX.  Also, the few data categories provided by Core are unique as they are all DCT (Data Category as Tag) style.  For example, instead of /term/ (a data category) being &lt;</w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:t>classificationElement</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:t xml:space="preserve"> type</w:t></w:r>' +
  '<w:proofErr w:type="gramStart"/>' +
  '<w:r><w:t>=&#8221;term</w:t></w:r>' +
  '<w:proofErr w:type="gramEnd"/>' +
  '<w:r><w:t>&#8221;&gt;, it is always &lt;term&gt;.</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
  '</w:p>'
Set-RangeXml $coreSub $coreBody

# ---------------------------------------------------------------------
# 4. Table cells: wrap the lone content run of each named paragraph in a
#    spellStart/spellEnd proofErr pair, keeping its original pPr intact.
# ---------------------------------------------------------------------

$rFontsRpr = '<w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:color w:val="000000"/></w:rPr>'

function Wrap-SimpleWordParagraph([int]$paraIndex, [string]$pPr, [string]$word) {
    $rng = $d.Paragraphs($paraIndex).Range
    $sub = $d.Range($rng.Start, $rng.End - 1)
    $body = '<w:p>' + $pPr + '<w:proofErr w:type="spellStart"/><w:r>' + $rFontsRpr + '<w:t>' + $word + '</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'
    Set-RangeXml $sub $body
}

$pPrA = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr><w:cnfStyle w:val="000000100000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="1" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/>' + $rFontsRpr + '</w:pPr>'
$pPrB = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr><w:cnfStyle w:val="000000000000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/>' + $rFontsRpr + '</w:pPr>'

# Paragraph 19: first "transacGrp"
Wrap-SimpleWordParagraph 19 $pPrA 'transacGrp'
# Paragraph 25: "conceptEntry"
Wrap-SimpleWordParagraph 25 $pPrB 'conceptEntry'
# Paragraph 26: "langSec"
Wrap-SimpleWordParagraph 26 $pPrB 'langSec'
# Paragraph 27: "termSec"
Wrap-SimpleWordParagraph 27 $pPrB 'termSec'
# Paragraph 29: "adminGrp" (28 "descripGrp" is left untouched)
Wrap-SimpleWordParagraph 29 $pPrB 'adminGrp'
# Paragraph 30: "termNoteGrp"
Wrap-SimpleWordParagraph 30 $pPrB 'termNoteGrp'
# Paragraph 31: "titleStmt"
Wrap-SimpleWordParagraph 31 $pPrB 'titleStmt'
# Paragraph 32: second "transacGrp"
Wrap-SimpleWordParagraph 32 $pPrB 'transacGrp'

# Paragraph 38: "    termSec" -- leading whitespace run kept, only the
# "termSec" run gets wrapped; different pPr (numId 1, hanging indent).
$pPrC = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:ind w:left="166" w:hanging="166"/><w:cnfStyle w:val="000000100000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="1" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/>' + $rFontsRpr + '</w:pPr>'
$rng38 = $d.Paragraphs(38).Range
$sub38 = $d.Range($rng38.Start, $rng38.End - 1)
$body38 = '<w:p>' + $pPrC + '<w:r>' + $rFontsRpr + '<w:t xml:space="preserve">    </w:t></w:r><w:proofErr w:type="spellStart"/><w:r>' + $rFontsRpr + '<w:t>termSec</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'
Set-RangeXml $sub38 $body38

Write-Host "Done"
